$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.056.60'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.783.74'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.54'
$ws.Range("E5").Value = '  +2.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.39'
$ws.Range("E8").Value = '  +4.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.290'
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0706'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("E11").Value = '  +1.75%  '
$ws.Range("D12").Value = '2.039.34'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.07'
$ws.Range("E13").Value = '  +5.97%  '
$ws.Range("D14").Value = '1.772.44'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.624'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '34.030.86'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.93'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.64'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").Value = '0.0₃0785'
$ws.Range("E20").Value = '  +2.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.75'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.34'
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0513'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.67'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.54'
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '1.397.18'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.661'
$ws.Range("E36").Value = '  +6.15%  '
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.26'
$ws.Range("E39").Value = '  +8.43%  '
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.916'
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '78.11'
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.45'
$ws.Range("E44").Value = '  +15.50%  '
$ws.Range("D45").Value = '0.0₆0140'
$ws.Range("E45").Value = '  +18.26%  '
$ws.Range("E46").Value = '  +3.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.30'
$ws.Range("E47").Value = '  +4.85%  '
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.86'
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").Value = '1.938.99'
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("E51").Value = '  +0.49%  '
